# Helper: build a string of $n repeated copies of character/string $ch
function Repeat($ch, $n) {
    $s = ""
    for ($i = 0; $i -lt $n; $i++) { $s = $s + $ch }
    return $s
}

# Helper: force-merge all runs spanning [start, end) of $d into a single run,
# while preserving the character formatting that already exists at `start`
# (the very first character of the span is left untouched, which anchors the
# formatting the engine re-uses when it rewrites the rest of the span).
function MergeRange($d, $start, $end) {
    $minEnd = $start + 1
    if ($end -le $minEnd) { return }
    $rest = $d.Range($start, $end)
    $restStart = $start + 1
    $rest = $d.Range($restStart, $end)
    $restOrig = $rest.Text
    $restLen = $restOrig.Length
    $placeholder = Repeat "Q" $restLen
    $rest.Text = $placeholder
    $placeholderLen = $placeholder.Length
    $restEnd = $restStart + $placeholderLen
    $rest2 = $d.Range($restStart, $restEnd)
    $rest2.Text = $restOrig
}

# Helper: split the run(s) covering position $pos into two pieces at $pos,
# without altering any text, by briefly planting and removing a bookmark.
function SplitAt($d, $pos) {
    $pt = $d.Range($pos, $pos)
    $pt.Bookmarks.Add("TempSplitMarker") | Out-Null
    $d.Bookmarks.Item("TempSplitMarker").Delete()
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Hyperlink text "https:/" + "/" + "mmarinov.netlify.com" (3 runs,
#    identical rPr) -> single run "https://mmarinov.netlify.com"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("https://mmarinov.netlify.com") | Out-Null
$h1start = $rng.Start
$h1end = $rng.End
MergeRange $d $h1start $h1end

# ---------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from inside
#    "Conducted individualized tutoring and s|et goals..." to inside
#    "Macaulay Honors College at Hunter College |<spaces>" (splitting
#    that whitespace run at the new bookmark location), AND merge the
#    "Conducted...and s" / "et goals..." runs back into one run.
# ---------------------------------------------------------------------
$goBackExists = $d.Bookmarks.Exists("_GoBack")
if ($goBackExists) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("Macaulay Honors College at Hunter College") | Out-Null
$afterCollege = $rng.End
$insertPos = $afterCollege + 1
$insertPt = $d.Range($insertPos, $insertPos)
$insertPt.Bookmarks.Add("_GoBack") | Out-Null

$rng = $d.Content
$rng.Find.Execute("Conducted individualized tutoring and s") | Out-Null
$condStart = $rng.Start
$rng2 = $d.Content
$rng2.Find.Execute("et goals for the chosen student") | Out-Null
$condEnd = $rng2.End
MergeRange $d $condStart $condEnd

# ---------------------------------------------------------------------
# 3) "xtracting keywords through " -> "xtracting keywords " + "through "
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("xtracting keywords ") | Out-Null
$splitPos3 = $rng.End
SplitAt $d $splitPos3

# ---------------------------------------------------------------------
# 4) "Created a" -> "Created " + "a"   (before " full stack app")
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Created a full stack") | Out-Null
$createdStart = $rng.Start
$splitPos4 = $createdStart + 8
SplitAt $d $splitPos4

# ---------------------------------------------------------------------
# 5) "Health, Weightlifting, Handball, Drawing, Ping Pong" ->
#    "Health, Weightlifting, Handball, Drawing" + ", Ping Pong"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Health, Weightlifting, Handball, Drawing") | Out-Null
$splitPos5 = $rng.End
SplitAt $d $splitPos5
